# Applies the Betfair odds update for Jogos_do_Dia_Betfair_Back_Lay_2026-01-15.xlsx
# Updates numeric odds/lay values for rows 2-12 and the Home/Away team names
# (plus kickoff Time) for the Friendly Matches rows that were reshuffled
# (rows 6-8), matching the new upstream feed snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.69
$ws.Range("G2").Value = 1.73
$ws.Range("H2").Value = 5.8
$ws.Range("I2").Value = 6.4
$ws.Range("L2").Value = 1.47
$ws.Range("N2").Value = 3.35
$ws.Range("O2").Value = 1.4
$ws.Range("P2").Value = 1.81
$ws.Range("Q2").Value = 2.18
$ws.Range("R2").Value = 1.29
$ws.Range("S2").Value = 4.2
$ws.Range("T2").Value = 2.08
$ws.Range("U2").Value = 1.83
$ws.Range("X2").Value = 12
$ws.Range("Z2").Value = 55
$ws.Range("AA2").Value = 190
$ws.Range("AC2").Value = 9
$ws.Range("AF2").Value = 9.199999999999999
$ws.Range("AH2").Value = 25
$ws.Range("AJ2").Value = 17.5
$ws.Range("AK2").Value = 20
$ws.Range("AM2").Value = 180
$ws.Range("AN2").Value = 13.5
$ws.Range("I3").Value = 13
$ws.Range("N3").Value = 9.4
$ws.Range("P3").Value = 3.9
$ws.Range("T3").Value = 1.7
$ws.Range("U3").Value = 2.22
$ws.Range("Y3").Value = 65
$ws.Range("AB3").Value = 17.5
$ws.Range("AH3").Value = 28
$ws.Range("AJ3").Value = 12
$ws.Range("AL3").Value = 29
$ws.Range("N4").Value = 1.28
$ws.Range("O4").Value = 1.16
$ws.Range("P4").Value = 1.28
$ws.Range("Q4").Value = 1.16
$ws.Range("S4").Value = 1.16
$ws.Range("N5").Value = 1.28
$ws.Range("P5").Value = 1.28
$ws.Range("D6").Value = 'Midtjylland'
$ws.Range("E6").Value = 'Ferencvaros'
$ws.Range("J6").Value = 1.09
$ws.Range("N6").Value = 1.1
$ws.Range("P6").Value = 2.88
$ws.Range("R6").Value = 1.76
$ws.Range("S6").Value = 1.05
$ws.Range("X6").Value = 970
$ws.Range("Y6").Value = 970
$ws.Range("AB6").Value = 970
$ws.Range("AC6").Value = 970
$ws.Range("AD6").Value = 970
$ws.Range("AG6").Value = 970
$ws.Range("AH6").Value = 970
$ws.Range("C7").Value = '11:00:00'
$ws.Range("D7").Value = 'Plzen'
$ws.Range("E7").Value = 'Sonderjyske'
$ws.Range("G7").Value = 600
$ws.Range("I7").Value = 870
$ws.Range("J7").Value = 1.04
$ws.Range("N7").Value = 1.25
$ws.Range("O7").Value = 1.14
$ws.Range("P7").Value = 1.24
$ws.Range("Q7").Value = 1.14
$ws.Range("S7").Value = 1.13
$ws.Range("C8").Value = '12:00:00'
$ws.Range("D8").Value = 'Puskas Akademia'
$ws.Range("E8").Value = 'Slovan Liberec'
$ws.Range("G8").Value = 970
$ws.Range("I8").Value = 970
$ws.Range("J8").Value = 1.09
$ws.Range("O8").Value = 1.18
$ws.Range("Q8").Value = 1.18
$ws.Range("S8").Value = 1.17
$ws.Range("X8").Value = 1000
$ws.Range("Y8").Value = 1000
$ws.Range("AB8").Value = 1000
$ws.Range("AC8").Value = 1000
$ws.Range("AD8").Value = 1000
$ws.Range("AG8").Value = 1000
$ws.Range("AH8").Value = 1000
$ws.Range("F9").Value = 3.8
$ws.Range("G9").Value = 3.85
$ws.Range("I9").Value = 2.32
$ws.Range("N9").Value = 2.98
$ws.Range("O9").Value = 1.49
$ws.Range("Q9").Value = 2.44
$ws.Range("V9").Value = 1.75
$ws.Range("W9").Value = 1.35
$ws.Range("X9").Value = 9.4
$ws.Range("AA9").Value = 30
$ws.Range("AE9").Value = 29
$ws.Range("AG9").Value = 16
$ws.Range("AO9").Value = 27
$ws.Range("H10").Value = 2.92
$ws.Range("I10").Value = 2.94
$ws.Range("U10").Value = 1.97
$ws.Range("X11").Value = 12
$ws.Range("L12").Value = 1.4
$ws.Range("M12").Value = 1.06
$ws.Range("N12").Value = 3.45
$ws.Range("P12").Value = 1.84
$ws.Range("R12").Value = 1.32
$ws.Range("S12").Value = 3.25
$ws.Range("X12").Value = 16
$ws.Range("Y12").Value = 980
$ws.Range("Z12").Value = 980
$ws.Range("AC12").Value = 10
$ws.Range("AD12").Value = 980
$ws.Range("AE12").Value = 100
$ws.Range("AF12").Value = 12
$ws.Range("AH12").Value = 980
$ws.Range("AK12").Value = 980
$ws.Range("AL12").Value = 980
$ws.Range("AN12").Value = 14
